$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, shifting existing rows 61..114 down to 62..115.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new data record.
$ws.Cells.Item(61, 1).Value = 7
$ws.Cells.Item(61, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(61, 3).Value = "Ñuble"
$ws.Cells.Item(61, 4).Value = 44907
$ws.Cells.Item(61, 5).Value = 16
$ws.Cells.Item(61, 6).Value = 100112021
$ws.Cells.Item(61, 7).Value = "Ají"
$ws.Cells.Item(61, 8).Value = "Americana (o)"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 100
$ws.Cells.Item(61, 11).Value = 15500
$ws.Cells.Item(61, 12).Value = 16000
$ws.Cells.Item(61, 13).Value = 15750
$ws.Cells.Item(61, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(61, 15).Value = "Región del Maule"
$ws.Cells.Item(61, 16).Value = 1050
$ws.Cells.Item(61, 17).Value = 15
$ws.Cells.Item(61, 18).Value = "Hortaliza"
